# Simulation_log.xlsx - "Add script to test successful breeding"
#
# Adds two new entries (rows 21 & 22) to the Script_tracking log describing
# scripts 022 (estimate lambda while accounting for skipped breeding) and
# 023 (neutral-lambda skipped-breeding model with failed breeders), renames
# the "TM" (Tailored model) abbreviation on the Notes sheet to "AM" (Adapted
# model), widens column D to fit the new dates, and leaves the
# Script_tracking sheet active/selected at the newly added row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Script_tracking")
$ws2 = $wb.Worksheets.Item("Notes")

# --- Notes sheet: rename abbreviation TM -> AM -------------------------
$ws2.Range("B5").Value = "AM"
$ws2.Range("C5").Value = "Adapted model. Means the CKMR model (with lambda parameter) was adapted to account for the process being tested with simulation (e.g. skipped-breeding)"
$ws2.Range("B6").Select() | Out-Null

# --- Script_tracking sheet: new row 21 (script 022) ---------------------
$ws1.Range("A21").Value = "Dovi IBS"
$ws1.Range("B21").Value = "Lemon Shark"
$ws1.Range("C21").Value = "Six"
$ws1.Range("D21").Value = 44419
$ws1.Range("E21").Value = "022_DoviIBS_LS_08.11.2021_SB_AM_estimate_lambda"
$ws1.Range("G21").Value = "This script attempts to estimate lambda along with abundance, while accounting for skipped breeding."
$ws1.Range("H21").Value = "Couldn't estimate lambda; it was confounded by the abundance estimate."
$ws1.Rows.Item(21).RowHeight = 30

# --- Script_tracking sheet: new row 22 (script 023) ---------------------
$ws1.Range("A22").Value = "Dovi IBS"
$ws1.Range("B22").Value = "Lemon Shark"
$ws1.Range("C22").Value = "Six"
$ws1.Range("D22").Value = 44480
$ws1.Range("E22").Value = "023_DoviIBS_LS_10.11.2021_neutral_lambda_SB_AM_failed_breeders"
$ws1.Range("G22").Value = "What happens when a certain percentage of females fail to breed? Will the model still work?"
$ws1.Rows.Item(22).RowHeight = 30

# --- Column D is now a bit wider to fit the new dates -------------------
$ws1.Columns.Item(4).AutoFit()

# --- Leave Script_tracking as the active sheet, selection on D22 --------
$ws1.Activate() | Out-Null
$ws1.Range("D22").Select() | Out-Null
